$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn Down Chart")

# ---------------------------------------------------------------------------
# 1. Update the title cell (Release number bumped from 2.0 to 1.0)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Project Title: Minesweeper`nRelease #: 1.0`nSprint #: 2"

# ---------------------------------------------------------------------------
# 2. Rewrite the Sprint Back Log rows (4-9) with the new M2 user stories
# ---------------------------------------------------------------------------
# Row 4 - M2-1
$ws.Range("B4").Value = "M2-1"
$ws.Range("C4").Value = "I would like to setup the database so we can save data"
$ws.Range("D4").Value = "Task 1"
$ws.Range("E4").Value = "Shawn"
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0

# Row 5 - M2-2
$ws.Range("B5").Value = "M2-2"
$ws.Range("C5").Value = "I would like to create login and registration controllers"
$ws.Range("D5").Value = "Task 2"
$ws.Range("E5").Value = "Shawn"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0

# Row 6 - M2-3
$ws.Range("B6").Value = "M2-3"
$ws.Range("C6").Value = "I would like to create models for user and registration"
$ws.Range("D6").Value = "Task 3"
$ws.Range("E6").Value = "Shawn"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0

# Row 7 - M2-4
$ws.Range("B7").Value = "M2-4"
$ws.Range("C7").Value = "I would like to create views that are consumed by controllers"
$ws.Range("D7").Value = "Task 4"
$ws.Range("E7").Value = "Richard"
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0

# Row 8 - M2-5 (previously blank)
$ws.Range("B8").Value = "M2-5"
$ws.Range("C8").Value = "I would like to insure that controllers and views validate data"
$ws.Range("D8").Value = "Task 5"
$ws.Range("E8").Value = "Richard"
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0

# Row 9 - M2-6 (previously blank)
$ws.Range("B9").Value = "M2-6"
$ws.Range("C9").Value = "I would like to update design documentation"
$ws.Range("D9").Value = "Task 6"
$ws.Range("E9").Value = "Richard"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0

# ---------------------------------------------------------------------------
# 3. Insert two new blank Sprint Back Log rows above the Estimate/Actual
#    summary rows (old rows 12/13 move down to 14/15), matching the
#    formatting used by the other blank rows (e.g. row 11).
# ---------------------------------------------------------------------------
$ws.Rows("12:13").Insert()

$ws.Range("B11:P11").Copy()
$ws.Range("B12:P12").PasteSpecial(-4122)
$ws.Range("B11:P11").Copy()
$ws.Range("B13:P13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Point the burn-down chart series at the relocated summary rows
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(""Ideal burndown"",,'Burn Down Chart'!`$G`$14:`$P`$14,1)"
$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(""Actual burndown"",,'Burn Down Chart'!`$G`$15:`$P`$15,2)"

# ---------------------------------------------------------------------------
# 5. Restore the active selection to C6 (was M6) and clear the scrolled
#    top-left cell so the sheet opens back at the top.
# ---------------------------------------------------------------------------
$ws.Range("C6").Select()
